$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "62.202.80"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.416.49"
$ws.Range("E3").Value = "  +1.30%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.13%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'562.47"
$ws.Range("E5").Value = "  +2.21%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'142.89"
$ws.Range("E6").Value = "  +1.08%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +1.63%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.411.27"
$ws.Range("E9").Value = "  +1.10%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.29%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -2.19%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  -0.31%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'25.66"
$ws.Range("E14").Value = "  -0.80%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +0.50%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").Value = "2.848.62"
$ws.Range("E16").Value = "  +1.12%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "62.153.47"
$ws.Range("E17").Value = "  +0.88%  "

# Row 18 - WrappedEther
$ws.Range("D18").Value = "2.411.81"
$ws.Range("E18").Value = "  +1.08%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'11.31"
$ws.Range("E19").Value = "  +2.05%  "

# Row 20 - swap BitcoinCash -> Polkadot
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.18"
$ws.Range("E20").Value = "  +0.53%  "

# Row 21 - swap Polkadot -> BitcoinCash
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'323.47"
$ws.Range("E21").Value = "  +0.44%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'6.84"
$ws.Range("E22").Value = "  +2.78%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.17%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'65.76"
$ws.Range("E24").Value = "  +2.07%  "

# Row 25 - SuiNetwork
$ws.Range("E25").Value = "  -1.77%  "

# Row 26 - Aptos
$ws.Range("D26").Value = "'9.01"
$ws.Range("E26").Value = "  -1.74%  "

# Row 27 - Bittensor
$ws.Range("D27").Value = "'577.51"
$ws.Range("E27").Value = "  +4.59%  "

# Row 28 - swap PEPE -> WrappedeETH
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.531.13"
$ws.Range("E28").Value = "  +2.32%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.17%  "

# Row 30 - swap WrappedeETH -> PEPE
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "0.0₃0948"
$ws.Range("E30").Value = "  +3.56%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'8.23"
$ws.Range("E31").Value = "  -0.40%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +1.53%  "

# Row 33 - Kaspa
$ws.Range("E33").Value = "  +0.62%  "

# Row 34 - PancakeSwap
$ws.Range("E34").Value = "  +0.89%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +0.32%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.13%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  -3.06%  "

# Row 38 - NEARProtocol
$ws.Range("D38").Value = "'4.72"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39 - PolygonEcosystemToken
$ws.Range("E39").Value = "  +0.04%  "

# Row 40 - Monero
$ws.Range("D40").Value = "'152.10"
$ws.Range("E40").Value = "  +4.68%  "

# Row 41 - EthereumClassic
$ws.Range("E41").Value = "  +0.74%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  -6.46%  "

# Row 43 - USDe
$ws.Range("E43").Value = "  -0.03%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "'2.30"
$ws.Range("E44").Value = "  +2.12%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'148.56"
$ws.Range("E45").Value = "  +0.27%  "

# Row 46 - Filecoin
$ws.Range("E46").Value = "  +0.82%  "

# Row 47 - Hedera
$ws.Range("D47").Value = "'0.0535"
$ws.Range("E47").Value = "  +1.12%  "

# Row 48 - InjectiveProtocol
$ws.Range("D48").Value = "'20.04"
$ws.Range("E48").Value = "  -0.45%  "

# Row 49 - Mantle
$ws.Range("E49").Value = "  +1.58%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +1.31%  "

# Row 51 - VeChain
$ws.Range("D51").Value = "'0.0228"
$ws.Range("E51").Value = "  +1.74%  "
